$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 92
$ws1.Range("F5").Value  = 1688
$ws1.Range("F6").Value  = 3274
$ws1.Range("F7").Value  = 885
$ws1.Range("F8").Value  = 2091
$ws1.Range("F9").Value  = 2004
$ws1.Range("F10").Value = 1035
$ws1.Range("F13").Value = 1625
$ws1.Range("F16").Value = 17
$ws1.Range("F17").Value = 75
$ws1.Range("F18").Value = 95
$ws1.Range("F19").Value = 1467
$ws1.Range("F20").Value = 544
$ws1.Range("F21").Value = 650
$ws1.Range("F23").Value = 11812
$ws1.Range("F24").Value = 11834
$ws1.Range("F27").Value = 1859
$ws1.Range("F29").Value = 468

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 65

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 65
$ws4.Range("F6").Value  = 92
$ws4.Range("F7").Value  = 1688
$ws4.Range("F8").Value  = 3274
$ws4.Range("F9").Value  = 885
$ws4.Range("F10").Value = 2091
$ws4.Range("F11").Value = 2004
$ws4.Range("F12").Value = 1035
$ws4.Range("F15").Value = 1625
$ws4.Range("F18").Value = 17
$ws4.Range("F20").Value = 75
$ws4.Range("F22").Value = 95
$ws4.Range("F23").Value = 1467
$ws4.Range("F24").Value = 544
$ws4.Range("F25").Value = 650
$ws4.Range("F27").Value = 11812
$ws4.Range("F28").Value = 11834
$ws4.Range("F31").Value = 1859
$ws4.Range("F35").Value = 468
